# Apply updated crypto price/volume data to sheet1 (Price=D, Volume(1h)=E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.709.38'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '2.208.35'
$ws.Range("E3").Value = '  -0.82%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = "'240.09"
$ws.Range("E5").Value = '  -1.88%  '

$ws.Range("D6").Value = "'0.619"
$ws.Range("E6").Value = '  -1.58%  '

$ws.Range("D7").Value = "'71.95"
$ws.Range("E7").Value = '  -2.35%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").Value = "'0.585"
$ws.Range("E9").Value = '  -4.76%  '

$ws.Range("D10").Value = "'41.00"
$ws.Range("E10").Value = '  -3.49%  '

$ws.Range("E11").Value = '  -3.28%  '

$ws.Range("E12").Value = '  +0.00%  '

$ws.Range("D13").Value = "'6.81"
$ws.Range("E13").Value = '  -4.39%  '

$ws.Range("D14").Value = '2.539.52'
$ws.Range("E14").Value = '  -0.73%  '

$ws.Range("D15").Value = "'14.02"
$ws.Range("E15").Value = '  -2.68%  '

$ws.Range("D16").Value = "'0.821"
$ws.Range("E16").Value = '  -3.50%  '

$ws.Range("D17").Value = '2.210.67'
$ws.Range("E17").Value = '  -1.18%  '

$ws.Range("D18").Value = '41.590.68'
$ws.Range("E18").Value = '  -1.27%  '

$ws.Range("E19").Value = '  -8.83%  '

$ws.Range("D20").Value = "'6.07"
$ws.Range("E20").Value = '  -1.27%  '

$ws.Range("D21").Value = "'71.21"
$ws.Range("E21").Value = '  -1.24%  '

$ws.Range("D22").Value = "'10.62"
$ws.Range("E22").Value = '  +7.04%  '

$ws.Range("D23").Value = "'227.22"
$ws.Range("E23").Value = '  -1.65%  '

$ws.Range("E24").Value = '  -6.54%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").Value = "'11.21"
$ws.Range("E26").Value = '  -5.50%  '

$ws.Range("E27").Value = '  -0.44%  '

$ws.Range("D28").Value = "'2.24"
$ws.Range("E28").Value = '  -2.37%  '

$ws.Range("E29").Value = '  -0.99%  '

$ws.Range("D30").Value = "'166.22"
$ws.Range("E30").Value = '  -0.46%  '

$ws.Range("D31").Value = "'20.24"
$ws.Range("E31").Value = '  -3.89%  '

$ws.Range("D32").Value = "'0.0784"
$ws.Range("E32").Value = '  -2.22%  '

$ws.Range("D33").Value = "'30.24"
$ws.Range("E33").Value = '  +2.69%  '

$ws.Range("D34").Value = "'5.29"
$ws.Range("E34").Value = '  -7.93%  '

$ws.Range("E35").Value = '  -1.94%  '

$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = '  -8.81%  '

$ws.Range("E37").Value = '  -4.65%  '

$ws.Range("E38").Value = '  -2.24%  '

$ws.Range("D39").Value = "'12.86"
$ws.Range("E39").Value = '  -1.08%  '

$ws.Range("E40").Value = '  -3.73%  '

$ws.Range("D41").Value = "'5.54"
$ws.Range("E41").Value = '  -1.34%  '

$ws.Range("D42").Value = "'63.02"
$ws.Range("E42").Value = '  +0.32%  '

$ws.Range("D43").Value = "'0.193"
$ws.Range("E43").Value = '  -3.75%  '

$ws.Range("D44").Value = "'8.56"
$ws.Range("E44").Value = '  -2.95%  '

$ws.Range("E45").Value = '  -2.67%  '

$ws.Range("D46").Value = "'100.59"
$ws.Range("E46").Value = '  -4.48%  '

$ws.Range("E47").Value = '  -1.54%  '

$ws.Range("E48").Value = '  -1.79%  '

$ws.Range("E49").Value = '  -3.53%  '

$ws.Range("E50").Value = '  -1.55%  '

$ws.Range("D51").Value = '2.417.77'
$ws.Range("E51").Value = '  -0.77%  '
